$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A19").Value = 43
$ws.Range("B19").Value = "small updates"
$ws.Range("C19").Value = "riya-morankar"
$ws.Range("D19").Value = "N/A"
$ws.Range("E19").Value = "edit1 to main"

# Force the date-looking string to be stored as text (not auto-converted
# to a date serial number) and keep the cell on the default "Normal"
# style, matching the other rows in this column.
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "2025-06-18"
$ws.Range("F19").Style = "Normal"
